$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- Title ---
Replace-Text "Nucleosynthesis: Birth of Stellar Atoms" "Journey into History: Unraveling the Chronicles of Human Endeavours"

# --- Author name: "Lucas Anderson" -> "Prof" / "." / " Jennifer Carter" (3 runs) ---
$pAuthor = $d.Paragraphs(2)
$authorStart = $pAuthor.Range.Start
$authorEnd = $pAuthor.Range.End
$rAuthor = $d.Range($authorStart, $authorEnd - 1)
$rAuthor.Text = "Prof"
$rDot = $d.Range($authorStart + 4, $authorStart + 4)
$rDot.InsertAfter(".")
$rRest = $d.Range($authorStart + 5, $authorStart + 5)
$rRest.InsertAfter(" Jennifer Carter")

# --- Email ---
Replace-Text "landerson@astralobservatory" "historyenlightenment101@educonnect"
$pEmail = $d.Paragraphs(3)
$emailFind = $pEmail.Range.Find
$null = $emailFind.Execute("org", $true, $false, $false, $false, $false, $true, 1, $false, "com", 2)

# --- Body paragraph (paragraph 5) sentence replacements (1:1, run count unchanged) ---
Replace-Text "Think about the atom that constitutes your body and the entirety of the physical world as we know it" "History, the profound chronicle of human experience, unravels the tapestry of our collective past and sheds light upon the present"
Replace-Text " Where do they come from? How did they come to exist, enabling the formation of stars, planets, and ultimately, ourselves? Nuclear fusion in stellar interiors provides the answer to these awe-inspiring questions" " It is a multifaceted panorama of civilizations, cultures, ideas, and events that have shaped the world we inhabit today"
Replace-Text " It is there, within the intense heat and pressure of stars, that simple atomic nuclei overcome their mutual repulsion and merge, fusing into heavier elements, facilitating the genesis of all elements beyond hydrogen and helium" " Embarking on this historical odyssey, we delve into the intricate narratives of influential individuals, transformative events, and enduring legacies"
Replace-Text " This captivating process, aptly named nucleosynthesis, played a pivotal role in shaping the elements of the universe, paving the way for the formation of the intricate structures observed in the cosmos" " Through the lens of historical inquiry, we decipher the enigma of humanity's triumphs and travails, exploring the intricacies of power, progress, and perseverance"
Replace-Text "Interwoven with the narrative of stellar nucleosynthesis is a tale of evolutionary cycles" "By delving into primary sources, we become immersed in the authentic voices of the past, gaining firsthand insights into the minds and motivations of historical actors"
Replace-Text " Massive stars, in their relentless pursuit of energy, fuse elements until they reach iron" " These documents, artifacts, and accounts provide tantalizing glimpses into bygone eras, enabling us to reconstruct the rich tapestry of human experience"
Replace-Text " This process subsequently ceases as no energy can be extracted from iron, leading to the implosion of the star in a cataclysmic spectacle known as a supernova" " Moreover, secondary sources, such as historical accounts, offer diverse perspectives and interpretations, inviting us to engage in critical analysis and fostering a deeper understanding of the past"
Replace-Text "Witnessing nucleosynthesis firsthand through astronomical observations is a privilege reserved for the most discerning instruments, capable of dissecting the light emitted by distant stars, unveiling the composition of their elemental tapestry" "Unveiling the complexities of historical phenomena, we uncover the intricate interplay of cause and effect"
Replace-Text " By meticulously scrutinizing the absorption and emission patterns within stellar spectra, astronomers can ascertain the chemical elements that reside within these celestial beacons, providing empirical evidence of the cosmic alchemy that has been ongoing for billions of years" " We scrutinize the factors that have shaped the course of history, delving into the social, political, economic, and cultural forces that have propelled or hindered human progress"

# Remove the leftover sentence (+ its leading space and trailing period) that
# has no counterpart in the new text: " It is during these explosive
# moments ... commencing the cycle anew."
$p5 = $d.Paragraphs(5)
$leftover = $d.Range($p5.Range.Start, $p5.Range.End)
$null = $leftover.Find.Execute(" It is during these explosive moments that the synthesized elements are expelled into the interstellar medium, enriching it with the lifeblood of heavy elements, ready to be taken up into subsequent generations of star formation, commencing the cycle anew.", $true)
$leftover.Text = ""

# The paragraph currently ends "...human progress." (that trailing "." is the
# pre-existing final-run period). Insert the new sentences *before* that
# final "." run, then re-add a "." to separate "human progress" from the
# first new sentence.
$p5 = $d.Paragraphs(5)
$insertPoint = $d.Range($p5.Range.End - 2, $p5.Range.End - 2)
$insertPoint.InsertAfter(".")
$p5 = $d.Paragraphs(5)
$insertPoint2 = $d.Range($p5.Range.End - 2, $p5.Range.End - 2)
$insertPoint2.InsertAfter(" This exploration illuminates the interconnectedness of events, revealing the ripple effects of decisions and actions, both grand and seemingly insignificant")
$p5 = $d.Paragraphs(5)
$insertPoint3 = $d.Range($p5.Range.End - 2, $p5.Range.End - 2)
$insertPoint3.InsertAfter(".")
$p5 = $d.Paragraphs(5)
$insertPoint4 = $d.Range($p5.Range.End - 2, $p5.Range.End - 2)
$insertPoint4.InsertAfter(" By examining past mistakes, we gain invaluable insights into the challenges and opportunities that lie ahead")

Write-Host "Stage 2 done"
Write-Host $d.Paragraphs(5).Range.Text
